# =======================================================================
# Instructions_left_hand.xlsx -- restructure table: reorder/add columns,
# add German localisation column, restyle header row, add new shared
# strings, resize columns, update selection.
# =======================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Text content blocks (literal here-strings -- no interpolation)
# ---------------------------------------------------------------------
$pic1 = @'
hljt_instr_images/instr_pic1.jpg
'@

$pic3 = @'
hljt_instr_images/instr_pic3.jpg
'@

$pic2_left = @'
hljt_instr_images/instr_pic2_left.jpg
'@

$en_instructions = @'
Instructions:

In this task, you will see images of left or right hands, their palms facing up or down. The images will be rotated at different angles.

Your task is to determine if the image corresponds to a left or right hand.

Your goal is to respond as quickly AND accurately as possible.

Each image will appear until you have responded. The next image will appear automatically.
'@

$en_hand = @'
Please use only your INDEX and MIDDLE FINGERS of the LEFT HAND to respond.

Place your index finger on the "H" key and the middle finger on the "G" key of your keyboard.

To respond:
Left Hand = G | H = Right Hand

You must maintain your hand on the keyboard throughout the task.

Keep your other hand on the desk, in the same position and as motionless as possible.
'@

$en_feedback = @'
After each image, you will receive a short feedback on your response:

If you respond correctly, the corresponding box will turn green

If you respond incorrectly, the corresponding box will turn red

Remember that your goal is to respond as accurately and fast as possible
'@

$es_instructions = @'
Instrucciones:

En esta tarea, verás imágenes del dorso o la palma de manos izquierdas o derechas. Las imágenes estarán rotadas con diferentes ángulos.

Tu tarea es determinar si la imagen corresponde a una mano izquierda o derecha.

Tu objetivo es responder lo más precisa y rápidamente posible.

Cada imagen aparecerá hasta que hayas respondido. Una vez hayas respondido, la siguiente imagen aparecerá automáticamente.
'@

$es_hand = @'
Por favor, usa solo los DEDOS ÍNDICE y CORAZÓN de tu MANO IZQUIERDA para responder.

Coloca tu dedo índice sobre la "H" y el corazón sobre la "G" de tu teclado.

Para responder:
Mano izquierda = G | H = Mano derecha

Tienes que mantener tu otra mano sobre la mesa durante la tarea.

Mantén tus manos en la misma posición y tan quietas como puedas.
'@

$es_feedback = @'
Tras cada imagen, recibirás un feedback corto sobre tu respuesta:

Si respondes correctamente, la caja correspondiente se volverá verde

Si respondes incorrectamente, la caja correspondiente se volverá roja

Recuerda que tu objetivo es responder tan precisa y rápidamente como puedas
'@

$fr_instructions = @'
Instructions :

Dans cette tâche, vous verrez des images de mains gauches ou droites, les paumes tournées vers le haut ou vers le bas. Les images seront tournées sous différents angles.

Votre tâche consiste à déterminer si l'image correspond à une main gauche ou droite.

Votre objectif est de répondre le plus rapidement ET le plus précisément possible.

Chaque image apparaîtra jusqu'à ce que vous ayez répondu. L'image suivante apparaîtra automatiquement.
'@

$fr_hand = @'
Veuillez utiliser uniquement L'INDEX et LE MAJEUR de votre MAIN GAUCHE pour répondre.

Placez votre index sur le "H" et votre majeur sur le "G" de votre clavier.

Pour répondre :
Main gauche = G | H = Main droite.

Vous devez garder votre autre main sur la table pendant la tâche.

Gardez vos mains dans la même position et aussi immobiles que possible.
'@

$fr_feedback = @'
Après chaque image, vous recevrez un bref commentaire sur votre réponse :

Si vous répondez correctement, la case correspondante devient verte.

Si votre réponse est incorrecte, la case correspondante devient rouge.

N'oubliez pas que votre objectif est de répondre le plus précisément et le plus rapidement possible.
'@

$de_instructions = @'
Anweisungen:
In dieser Aufgabe sehen Sie Bilder von linken oder rechten Händen, deren Handflächen nach oben oder unten zeigen. Die Bilder werden in verschiedenen Winkeln gedreht.
Ihre Aufgabe ist es zu bestimmen, ob das Bild einer linken oder rechten Hand entspricht.
Ihr Ziel ist es, sowohl SCHNELL als auch GENAU zu antworten.
Jedes Bild wird angezeigt, bis Sie geantwortet haben. Das nächste Bild erscheint automatisch.
'@

$de_hand = @'
Bitte verwenden Sie nur Ihren ZEIGEFINGER und MITTELFINGER der LINKEN HAND, um zu antworten.
Legen Sie Ihren Zeigefinger auf die Taste „H“ und den Mittelfinger auf die Taste „G“ Ihrer Tastatur.
Zum Antworten:
Linke Hand = G | H = Rechte Hand
Sie müssen Ihre Hand während der gesamten Aufgabe auf der Tastatur halten.
Halten Sie Ihre andere Hand auf dem Tisch, in derselben Position und so ruhig wie möglich.
'@

$de_feedback = @'
Nach jedem Bild erhalten Sie ein kurzes Feedback zu Ihrer Antwort:
Wenn Sie korrekt antworten, wird das entsprechende Feld grün.
Wenn Sie falsch antworten, wird das entsprechende Feld rot.
Denken Sie daran, dass Ihr Ziel darin besteht, so genau und schnell wie möglich zu antworten.
'@

# ---------------------------------------------------------------------
# 2) Header labels (row 1) -- new column order: pics, EN, ES, FR, DE, w, h
#    (values are entered in this particular order so that newly-created
#    shared-string table entries line up with the canonical file)
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "inst_msg_EN"
$ws.Range("C1").Value = "inst_msg_ES"
$ws.Range("D1").Value = "inst_msg_FR"
$ws.Range("F1").Value = "image_w"
$ws.Range("G1").Value = "image_h"
$ws.Range("A1").Value = "inst_pics"
$ws.Range("E1").Value = "inst_msg_DE"

# ---------------------------------------------------------------------
# 3) Data rows 2-4
# ---------------------------------------------------------------------
$ws.Range("A2").Value = $pic1
$ws.Range("B2").Value = $en_instructions
$ws.Range("C2").Value = $es_instructions
$ws.Range("D2").Value = $fr_instructions
$ws.Range("E2").Value = $de_instructions
$ws.Range("F2").Value = 0.6
$ws.Range("G2").Value = 0.5

$ws.Range("A3").Value = $pic2_left
$ws.Range("B3").Value = $en_hand
$ws.Range("C3").Value = $es_hand
$ws.Range("D3").Value = $fr_hand
$ws.Range("E3").Value = $de_hand
$ws.Range("F3").Value = 0.6
$ws.Range("G3").Value = 0.5

$ws.Range("A4").Value = $pic3
$ws.Range("B4").Value = $en_feedback
$ws.Range("C4").Value = $es_feedback
$ws.Range("D4").Value = $fr_feedback
$ws.Range("E4").Value = $de_feedback
$ws.Range("F4").Value = 0.6
$ws.Range("G4").Value = 0.5

# ---------------------------------------------------------------------
# 4) Rich-text colour runs inside the feedback / hand messages
#    (these exact runs are carried over unmodified from the source file)
# ---------------------------------------------------------------------
# B4 = EN feedback: "green" in green, "red" in red
$ws.Range("B4").Characters(129, 8).Font.Color = 5287936
$ws.Range("B4").Characters(197, 4).Font.Color = 255

# C4 = ES feedback: "verde" in green, "roja" in red
$ws.Range("C4").Characters(131, 5).Font.Color = 5287936
$ws.Range("C4").Characters(203, 4).Font.Color = 255

# D4 = FR feedback: "verte" in green, "rouge" in red
$ws.Range("D4").Characters(138, 5).Font.Color = 5287936
$ws.Range("D4").Characters(210, 5).Font.Color = 255

# B3 = EN hand message: "INDEX and MIDDLE FINGERS of the LEFT HAND" in bold
$ws.Range("B3").Characters(22, 41).Font.Bold = $true

# ---------------------------------------------------------------------
# 5) Cell styles: wrap text on the 4 language columns, plain on the
#    picture-filename column and the numeric size columns
# ---------------------------------------------------------------------
foreach ($r in 2..4) {
    $ws.Range("A$r").Style = "Normal"
    foreach ($col in @("B","C","D","E")) {
        $ws.Range("$col$r").WrapText = $true
    }
    $ws.Range("F$r").Style = "Normal"
    $ws.Range("G$r").Style = "Normal"
}

# Re-assert the fixed row heights (typing the long strings above would
# otherwise have triggered Excels automatic row auto-fit)
$ws.Rows.Item(2).RowHeight = 24.5
$ws.Rows.Item(3).RowHeight = 24.5
$ws.Rows.Item(4).RowHeight = 24.5

# ---------------------------------------------------------------------
# 6) Header row style: solid blue fill + white font
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:G1")
$headerRange.Interior.Color = 15773696
$headerRange.Font.ThemeColor = 2

# ---------------------------------------------------------------------
# 7) Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 33.90625
$ws.Columns.Item(2).ColumnWidth = 16.81640625
$ws.Columns.Item(3).ColumnWidth = 15.7265625
$ws.Columns.Item(4).ColumnWidth = 14.6328125
$ws.Columns.Item(5).ColumnWidth = 14.6328125
$ws.Columns.Item(6).ColumnWidth = 7.436197916666667
$ws.Columns.Item(7).ColumnWidth = 7.072916666666667

# ---------------------------------------------------------------------
# 8) Selection / active cell
# ---------------------------------------------------------------------
$ws.Range("E4").Select()
